# cancermine v49 + custom aliases
# - Bump cancermine version string on the "basic" sheet (row 2, column E)
# - Bump uniprot version string on the "gencode" sheet (row 4, column E)
# - Update selection on "basic" sheet to E2
# - Move the active worksheet / selected tab from "predisposition" to "gencode"

$wb = $excel.ActiveWorkbook

# --- basic sheet: cancermine version bump + selection ---
$basic = $wb.Worksheets.Item("basic")
$basic.Range("E2").Value = "v49 (January 2023)"
$basic.Range("E2").Select()

# --- gencode sheet: uniprot version bump ---
$gencode = $wb.Worksheets.Item("gencode")
$gencode.Range("E4").Value = "2022_05"

# --- make "gencode" the active sheet/tab (was "predisposition") ---
$gencode.Activate()
